# Regen save_data: recompute column G ("K") values for rows 2-61
# (commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(1,0,3,0,1,1,1,1,1,2,2,2,0,1,1,0,1,1,0,0,3,1,2,0,1,2,1,0,1,1,4,1,2,2,5,2,0,1,1,1,0,1,2,1,0,1,3,0,1,1,1,1,0,6,4,1,2,0,1,1)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
